$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.201154884288428
$ws.Range("C2").Value = 0.3177404093272571
$ws.Range("D2").Value = 0.02876610486227094
$ws.Range("F2").Value = 0.3236978582895986
$ws.Range("G2").Value = 0.002372653629534434
$ws.Range("I2").Value = 0.2282277210525693
$ws.Range("O2").Value = 0.9993314048698494

$ws.Range("B3").Value = 1.051038546440168
$ws.Range("C3").Value = 0.2789609095972594
$ws.Range("D3").Value = 0.02518192536900443
$ws.Range("F3").Value = 0.3273827086035652
$ws.Range("G3").Value = 0.002375057464057778
$ws.Range("I3").Value = 0.2361230433375314
$ws.Range("O3").Value = 1.021985576090842

$ws.Range("B4").Value = 0.9585017898886576
$ws.Range("C4").Value = 0.2550418583422527
$ws.Range("D4").Value = 0.02297188721951215
$ws.Range("F4").Value = 0.3300656853817081
$ws.Range("G4").Value = 0.002376610789542919
$ws.Range("I4").Value = 0.2412824757097596
$ws.Range("O4").Value = 1.037261736567771

$ws.Range("B5").Value = 0.9207033607141284
$ws.Range("C5").Value = 0.2452681639283583
$ws.Range("D5").Value = 0.02206899705085164
$ws.Range("F5").Value = 0.3312643757020162
$ws.Range("G5").Value = 0.002377263291422302
$ws.Range("I5").Value = 0.2434630977389847
$ws.Range("O5").Value = 1.043829181787331

$ws.Range("B6").Value = 0.9144216697878846
$ws.Range("C6").Value = 0.243643671622209
$ws.Range("D6").Value = 0.02191893697894898
$ws.Range("F6").Value = 0.3314697695418438
$ws.Range("G6").Value = 0.002377372818615185
$ws.Range("I6").Value = 0.2438298984441293
$ws.Range("O6").Value = 1.044940335818154

$ws.Range("B7").Value = 0.9579923835280511
$ws.Range("C7").Value = 0.2549101531377573
$ws.Range("D7").Value = 0.02295971967703281
$ws.Range("F7").Value = 0.3300814252546687
$ws.Range("G7").Value = 0.002376619510314891
$ws.Range("I7").Value = 0.2413115684171009
$ws.Range("O7").Value = 1.037348923267047

$ws.Range("B8").Value = 1.149472072274534
$ws.Range("C8").Value = 0.3043921429664067
$ws.Range("D8").Value = 0.02753225631741429
$ws.Range("F8").Value = 0.3248809030367141
$ws.Range("G8").Value = 0.002373466448931939
$ws.Range("I8").Value = 0.2308851846948075
$ws.Range("O8").Value = 1.006858181863279

$ws.Range("B9").Value = 1.521966455510494
$ws.Range("C9").Value = 0.4005391612835183
$ws.Range("D9").Value = 0.03642241861982143
$ws.Range("F9").Value = 0.3180358742550169
$ws.Range("G9").Value = 0.002367894560908106
$ws.Range("I9").Value = 0.2129233983883974
$ws.Range("O9").Value = 0.9579659265190799

$ws.Range("B10").Value = 1.793700222784082
$ws.Range("C10").Value = 0.4706071551510149
$ws.Range("D10").Value = 0.04290464719974807
$ws.Range("F10").Value = 0.3150747288718136
$ws.Range("G10").Value = 0.002364169913365688
$ws.Range("I10").Value = 0.2012564359015838
$ws.Range("O10").Value = 0.9287689650883379

$ws.Range("B11").Value = 1.916876674765831
$ws.Range("C11").Value = 0.5023529229084716
$ws.Range("D11").Value = 0.04584233398222182
$ws.Range("F11").Value = 0.3141817210817237
$ws.Range("G11").Value = 0.002362554832210198
$ws.Range("I11").Value = 0.1962841619127689
$ws.Range("O11").Value = 0.9169634728299627

$ws.Range("B12").Value = 1.96345525837603
$ws.Range("C12").Value = 0.5143551002347522
$ws.Range("D12").Value = 0.04695310438441425
$ws.Range("F12").Value = 0.3139092530941809
$ws.Range("G12").Value = 0.002361954585517519
$ws.Range("I12").Value = 0.1944497522396482
$ws.Range("O12").Value = 0.9127067239338515

$ws.Range("B13").Value = 1.953426693085248
$ws.Range("C13").Value = 0.5117710852893538
$ws.Range("D13").Value = 0.04671395538950662
$ws.Range("F13").Value = 0.3139650057378987
$ws.Range("G13").Value = 0.002362083355277096
$ws.Range("I13").Value = 0.1948426641743147
$ws.Range("O13").Value = 0.9136139622536632

$ws.Range("B14").Value = 1.920710055318011
$ws.Range("C14").Value = 0.5033407395944778
$ws.Range("D14").Value = 0.04593375155732815
$ws.Range("F14").Value = 0.3141579860454797
$ws.Range("G14").Value = 0.002362505222351285
$ws.Range("I14").Value = 0.1961322705586134
$ws.Range("O14").Value = 0.9166089757667777

$ws.Range("B15").Value = 1.900661533356356
$ws.Range("C15").Value = 0.4981743788421795
$ws.Range("D15").Value = 0.04545563491252835
$ws.Range("F15").Value = 0.3142847594245168
$ws.Range("G15").Value = 0.002362765104751185
$ws.Range("I15").Value = 0.1969285143941207
$ws.Range("O15").Value = 0.9184713851303741

$ws.Range("B16").Value = 1.785641319351839
$ws.Range("C16").Value = 0.4685298420550339
$ws.Range("D16").Value = 0.04271243287098514
$ws.Range("F16").Value = 0.3151422600419949
$ws.Range("G16").Value = 0.002364277053571623
$ws.Range("I16").Value = 0.2015881539367984
$ws.Range("O16").Value = 0.9295703127746577

$ws.Range("B17").Value = 1.714966259300013
$ws.Range("C17").Value = 0.4503104090934471
$ws.Range("D17").Value = 0.04102667136199045
$ws.Range("F17").Value = 0.3157848868165658
$ws.Range("G17").Value = 0.00236522485427379
$ws.Range("I17").Value = 0.2045327272972621
$ws.Range("O17").Value = 0.9367582995229071

$ws.Range("B18").Value = 1.674274909528322
$ws.Range("C18").Value = 0.4398190227020109
$ws.Range("D18").Value = 0.04005602421131016
$ws.Range("F18").Value = 0.3161972203430778
$ws.Range("G18").Value = 0.002365777469422182
$ws.Range("I18").Value = 0.2062578901042973
$ws.Range("O18").Value = 0.9410315151751121

$ws.Range("B19").Value = 1.660490577145993
$ws.Range("C19").Value = 0.4362647736163581
$ws.Range("D19").Value = 0.03972720298990851
$ws.Range("F19").Value = 0.3163441527479947
$ws.Range("G19").Value = 0.002365965859376747
$ws.Range("I19").Value = 0.2068474054620574
$ws.Range("O19").Value = 0.9425021607101201

$ws.Range("B20").Value = 1.722493996600065
$ws.Range("C20").Value = 0.4522511516556165
$ws.Range("D20").Value = 0.04120623187525041
$ws.Range("F20").Value = 0.3157120546626544
$ws.Range("G20").Value = 0.00236512318704306
$ws.Range("I20").Value = 0.2042160078603359
$ws.Range("O20").Value = 0.9359787425041191

$ws.Range("B21").Value = 1.930321529869047
$ws.Range("C21").Value = 0.5058174652274943
$ws.Range("D21").Value = 0.04616296215409932
$ws.Range("F21").Value = 0.3140995169001997
$ws.Range("G21").Value = 0.002362381002078568
$ws.Range("I21").Value = 0.1957521637299739
$ws.Range("O21").Value = 0.9157234555799647

$ws.Range("B22").Value = 2.065764770459225
$ws.Range("C22").Value = 0.5407135543742356
$ws.Range("D22").Value = 0.04939271845672977
$ws.Range("F22").Value = 0.3134287351737299
$ws.Range("G22").Value = 0.00236065495525013
$ws.Range("I22").Value = 0.1905032836665822
$ws.Range("O22").Value = 0.9037319263363202

$ws.Range("B23").Value = 1.993512251198297
$ws.Range("C23").Value = 0.5220994116104407
$ws.Range("D23").Value = 0.04766985216707553
$ws.Range("F23").Value = 0.3137515566942639
$ws.Range("G23").Value = 0.002361570144877598
$ws.Range("I23").Value = 0.1932787415729224
$ws.Range("O23").Value = 0.9100175060132187

$ws.Range("B24").Value = 1.719090891113126
$ws.Range("C24").Value = 0.4513737941924774
$ws.Range("D24").Value = 0.04112505717365877
$ws.Range("F24").Value = 0.3157448485464585
$ws.Range("G24").Value = 0.00236516912665691
$ws.Range("I24").Value = 0.2043590961642585
$ws.Range("O24").Value = 0.9363307419269091

$ws.Range("B25").Value = 1.421529257859163
$ws.Range("C25").Value = 0.3746269847599706
$ws.Range("D25").Value = 0.03402587212623587
$ws.Range("F25").Value = 0.3195261938634602
$ws.Range("G25").Value = 0.002369336840869935
$ws.Range("I25").Value = 0.2175150168982812
$ws.Range("O25").Value = 0.9700174400337289
